$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("D3").Value = -10
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 20

# Update the active selection on the sheet
$ws.Range("B7").Select()
